$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure columns C:G keep their text (shared-string) representation instead of
# being auto-converted to numbers/percentages when we assign the new values.
$dataRange = $ws.Range("C2:G13")
$dataRange.NumberFormat = "@"

# Update team names for rows 5 and 6 (A. Klagenfurt / Hartberg swapped order).
$ws.Range("B5").Value = "A. Klagenfurt"
$ws.Range("B6").Value = "Hartberg"

# Update the stats columns (Cartoes, Escanteios, 1.5+, 2.5+, Med. Gols) with
# the refreshed 15-03-2024 values.
$ws.Range("D2").Value = "5.0"
$ws.Range("E2").Value = "72%"
$ws.Range("F2").Value = "43%"
$ws.Range("G2").Value = "2.67"
$ws.Range("C3").Value = "2.1"
$ws.Range("E3").Value = "77%"
$ws.Range("G3").Value = "2.38"
$ws.Range("C4").Value = "2.4"
$ws.Range("D4").Value = "5.5"
$ws.Range("E4").Value = "62%"
$ws.Range("F4").Value = "28%"
$ws.Range("G4").Value = "2.05"
$ws.Range("C5").Value = "2.0"
$ws.Range("D5").Value = "4.6"
$ws.Range("E5").Value = "67%"
$ws.Range("F5").Value = "52%"
$ws.Range("G5").Value = "2.57"
$ws.Range("C6").Value = "2.4"
$ws.Range("D6").Value = "4.6"
$ws.Range("E6").Value = "71%"
$ws.Range("F6").Value = "62%"
$ws.Range("G6").Value = "2.81"
$ws.Range("C7").Value = "2.2"
$ws.Range("E7").Value = "76%"
$ws.Range("F7").Value = "38%"
$ws.Range("G7").Value = "2.71"
$ws.Range("C8").Value = "3.0"
$ws.Range("D8").Value = "5.7"
$ws.Range("E8").Value = "66%"
$ws.Range("F8").Value = "48%"
$ws.Range("G8").Value = "2.14"
$ws.Range("C9").Value = "2.6"
$ws.Range("D9").Value = "5.3"
$ws.Range("E9").Value = "77%"
$ws.Range("F9").Value = "58%"
$ws.Range("G9").Value = "2.81"
$ws.Range("C10").Value = "2.0"
$ws.Range("D10").Value = "3.6"
$ws.Range("E10").Value = "81%"
$ws.Range("G10").Value = "2.76"
$ws.Range("C11").Value = "2.8"
$ws.Range("D11").Value = "5.2"
$ws.Range("E11").Value = "67%"
$ws.Range("F11").Value = "43%"
$ws.Range("G11").Value = "2.14"
$ws.Range("C12").Value = "2.3"
$ws.Range("D12").Value = "4.4"
$ws.Range("E12").Value = "75%"
$ws.Range("F12").Value = "48%"
$ws.Range("G12").Value = "2.86"
$ws.Range("C13").Value = "2.4"
$ws.Range("D13").Value = "2.9"
$ws.Range("E13").Value = "76%"
$ws.Range("F13").Value = "49%"
$ws.Range("G13").Value = "2.86"

# Restore the default (General) style now that the values have been written,
# so the cells end up without an explicit text number format applied.
$dataRange.Style = "Normal"
